# Updates cryptos list data: prices (D) and 1h volume-change percentages (E)
# for several rows, plus a rank swap between RenderToken/Fetch.AI (rows 27-28)
# and between FLOKI/Mantle (rows 43-44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.954.64'
$ws.Range("E2").Value = '  +3.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.803.81'
$ws.Range("E3").Value = '  +1.30%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '697.44'
$ws.Range("E5").Value = '  +11.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.10'
$ws.Range("E6").Value = '  +5.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.802.89'
$ws.Range("E7").Value = '  +1.43%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +1.77%  '

$ws.Range("E10").Value = '  +3.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.44'
$ws.Range("E11").Value = '  +8.05%  '

$ws.Range("E12").Value = '  +1.79%  '

$ws.Range("E13").Value = '  +10.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.41'
$ws.Range("E14").Value = '  +5.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.446.63'
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.803.40'
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.927.20'
$ws.Range("E17").Value = '  +3.20%  '

$ws.Range("E18").Value = '  +1.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.23'
$ws.Range("E19").Value = '  +3.82%  '

$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("E21").Value = '  +17.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.07'
$ws.Range("E22").Value = '  +3.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("E23").Value = '  +2.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.02'
$ws.Range("E24").Value = '  +3.01%  '

$ws.Range("E25").Value = '  +3.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.43'
$ws.Range("E26").Value = '  +3.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.954.73'
$ws.Range("E29").Value = '  +1.15%  '

$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  +14.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.55'
$ws.Range("E32").Value = '  +6.79%  '

$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.65'
$ws.Range("E34").Value = '  +4.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.181'
$ws.Range("E35").Value = '  +2.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.25'
$ws.Range("E36").Value = '  +4.96%  '

$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.754.37'
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("E39").Value = '  +3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.53'
$ws.Range("E40").Value = '  +9.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.99'
$ws.Range("E41").Value = '  +4.47%  '

$ws.Range("E42").Value = '  +13.16%  '

$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '163.54'
$ws.Range("E47").Value = '  +5.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.27'
$ws.Range("E48").Value = '  +5.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.04'
$ws.Range("E49").Value = '  +2.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.301'
$ws.Range("E50").Value = '  +3.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.39'
$ws.Range("E51").Value = '  +0.05%  '

# Row 27 and 28: RenderToken/Fetch.AI swap positions with updated values
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.50'
$ws.Range("E27").Value = '  +4.58%  '

$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.17'
$ws.Range("E28").Value = '  +3.39%  '

# Row 43 and 44: Mantle/FLOKI swap positions with updated values
$ws.Range("B43").Value = 'FLOKI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000330'
$ws.Range("E43").Value = '  +24.75%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.974'
$ws.Range("E44").Value = '  +2.04%  '

